$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
if (-not $ws) { $ws = $wb.ActiveSheet }

# --- Simple value updates (no structural shift) ---
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-27T12:23:18-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Insert a new "Jurisdiction" row after the Contact rows (new row 12),
#     pushing Description / Purpose / Copyright / Immutable down by one. ---
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Re-apply the same formatting used by the surrounding data rows, since a
# freshly inserted row does not automatically inherit it.
$ws.Range("A13:B13").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)
